$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Price (D) and Volume(1h) (E) columns ---
# Plain-decimal-looking Price strings are written with a leading apostrophe
# (Excel text-literal marker) so they stay text instead of becoming numbers,
# then the style is reset to Normal so no stray number-format style is left behind.
$ws.Range("D2").Value = "26.487.79"
$ws.Range("E2").Value = "  -0.70%  "

$ws.Range("D3").Value = "1.626.20"
$ws.Range("E3").Value = "  -0.54%  "

$ws.Range("E4").Value = "  +0.11%  "

$ws.Range("E5").Value = "  +0.02%  "

$ws.Range("E6").Value = "  +1.35%  "

$ws.Range("E7").Value = "  +0.07%  "

$ws.Range("E8").Value = "  -0.10%  "

$ws.Range("E9").Value = "  -1.75%  "

$ws.Range("D10").Value = "'18.75"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.41%  "

$ws.Range("E11").Value = "  +0.50%  "

$ws.Range("D12").Value = "1.851.25"
$ws.Range("E12").Value = "  -0.64%  "

$ws.Range("D13").Value = "1.648.25"
$ws.Range("E13").Value = "  +1.19%  "

$ws.Range("D14").Value = "'4.13"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.52%  "

$ws.Range("D15").Value = "'0.522"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.53%  "

$ws.Range("D16").Value = "'64.92"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.07%  "

$ws.Range("D17").Value = "26.491.60"
$ws.Range("E17").Value = "  -0.69%  "

$ws.Range("D18").Value = "0.0₃0739"

$ws.Range("D19").Value = "'214.36"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.92%  "

$ws.Range("E20").Value = "  +0.19%  "

$ws.Range("E21").Value = "  -0.54%  "

$ws.Range("D22").Value = "'6.24"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.28%  "

$ws.Range("E23").Value = "  -1.17%  "

$ws.Range("D24").Value = "'2.05"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +7.06%  "

$ws.Range("D25").Value = "'148.59"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.43%  "

$ws.Range("E26").Value = "  +0.19%  "

$ws.Range("E27").Value = "  -0.37%  "

$ws.Range("D28").Value = "'6.85"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.59%  "

$ws.Range("D29").Value = "'15.52"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.87%  "

$ws.Range("D30").Value = "'0.0508"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.76%  "

$ws.Range("E31").Value = "  -1.08%  "

$ws.Range("D32").Value = "'3.33"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.01%  "

$ws.Range("D33").Value = "'2.95"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.44%  "

$ws.Range("D34").Value = "1.237.54"
$ws.Range("E34").Value = "  +5.97%  "

$ws.Range("D35").Value = "'1.50"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.00%  "

$ws.Range("E36").Value = "  -1.84%  "

$ws.Range("E37").Value = "  +3.91%  "

$ws.Range("E38").Value = "  +0.12%  "

$ws.Range("D39").Value = "'0.507"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.76%  "

$ws.Range("D40").Value = "'0.793"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.56%  "

$ws.Range("E41").Value = "  -1.88%  "

$ws.Range("D42").Value = "'0.795"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.19%  "

$ws.Range("E43").Value = "  -0.83%  "

$ws.Range("D44").Value = "1.761.03"
$ws.Range("E44").Value = "  -0.66%  "

$ws.Range("D45").Value = "'92.91"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.56%  "

$ws.Range("E46").Value = "  +2.16%  "

# --- Row 47/48 swap: Aave and BabyDogeCoin swap positions, with updated D/E values ---
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").Value = "0.0₆0104"
$ws.Range("E47").Value = "  -0.02%  "

$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").Value = "'54.80"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.13%  "

$ws.Range("D49").Value = "'0.0508"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.68%  "

$ws.Range("D50").Value = "'0.407"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.64%  "

$ws.Range("D51").Value = "'7.45"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.30%  "
